# Replaced buck converter inductor
# Update the efficiency values for the buck converter rows to reflect the
# new inductor, then move the selection to L5 as a result of editing L-column
# cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buck converters")
$ws.Activate()

# Efficiency (column L) updates for rows 3 and 4
$ws.Range("L3").Value = 0.84
$ws.Range("L4").Value = 0.78

# Leave the selection on L5, matching the post-edit cursor position
$ws.Range("L5").Select()
